$d = $word.ActiveDocument

# 1. Remove the hidden "_GoBack" bookmark (it is name-addressable even
#    though it does not show up in Bookmarks.Count / iteration).
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# 2. Split the title run "Задание к модулю " (chars 0..16 of paragraph 1)
#    into six separate runs - "Задание" / " " / "к" / " " / "модулю" / " " -
#    each carrying the same <w:lang w:val="ru-RU"/> run properties, while
#    keeping the following "ADO.Net." run intact. We target the whole
#    "Задание к модулю ADO.Net." span (0..25) and re-supply "ADO.Net." verbatim
#    in the replacement markup, because InsertXML on a range whose Start sits
#    exactly on a run boundary relocates the inserted content to the end of
#    that run instead of the deletion point - including both runs in one
#    InsertXML call keeps everything contiguous and in the right order.
$r = $d.Range(0, 25)
$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t>Задание</w:t></w:r><w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t>к</w:t></w:r><w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t>модулю</w:t></w:r><w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t>ADO.Net.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r.InsertXML($xml)
